$d = $word.ActiveDocument

# The first two paragraphs of the document are the pandoc-style title
# block: a Heading1 "Fall Appeal - 1973" (wrapped in a bookmark) and a
# bold "By Dorothy Day" byline. We convert them into a real pandoc title
# block: a Title-styled paragraph with just the title text, and a
# separate Authors-styled paragraph with just the author's name (the
# "By " prefix is dropped). Each paragraph's text is split word-by-word
# into its own runs (with separate run(s) for the intervening spaces/
# hyphen), matching how pandoc emits these titles.

$p1 = $d.Paragraphs.Item(1)
$p2 = $d.Paragraphs.Item(2)

$startPos = $p1.Range.Start
$endPos = $p2.Range.End
$rng = $d.Range($startPos, $endPos)

$xmlFrag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t xml:space="preserve">Fall</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Appeal</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">-</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">1973</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Authors"/></w:pPr><w:r><w:t xml:space="preserve">Dorothy</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Day</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$rng.InsertXML($xmlFrag)
